$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.493.69"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.381.17"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'552.27"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").Value = "'140.29"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.525"
$ws.Range("E8").Value = "  +2.73%  "
$ws.Range("D9").Value = "2.381.33"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").Value = "'0.108"
$ws.Range("E10").Value = "  +6.55%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "'5.35"
$ws.Range("E12").Value = "  +3.70%  "
$ws.Range("E13").Value = "  +5.28%  "
$ws.Range("D14").Value = "'25.58"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("E15").Value = "  +7.23%  "
$ws.Range("D16").Value = "61.388.23"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "'10.99"
$ws.Range("E17").Value = "  +5.72%  "
# Rows 18-19: BitcoinCash and Polkadot swapped places
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'4.15"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'321.23"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").Value = "'6.77"
$ws.Range("E20").Value = "  +5.40%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'1.77"
$ws.Range("E22").Value = "  -4.49%  "
$ws.Range("D23").Value = "'64.37"
$ws.Range("E23").Value = "  +3.04%  "
$ws.Range("D24").Value = "'8.94"
$ws.Range("E24").Value = "  +12.66%  "
$ws.Range("D25").Value = "'8.22"
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("D26").Value = "'521.63"
$ws.Range("E26").Value = "  +4.96%  "
$ws.Range("D27").Value = "0.0₃0905"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("E28").Value = "  +6.19%  "
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("E30").Value = "  +4.33%  "
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'5.58"
$ws.Range("E33").Value = "  +8.08%  "
$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  +6.56%  "
$ws.Range("E35").Value = "  +8.89%  "
$ws.Range("D36").Value = "'0.379"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "'18.56"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "'146.60"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'41.35"
$ws.Range("D41").Value = "'148.54"
$ws.Range("E41").Value = "  +9.60%  "
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  +7.26%  "
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("D44").Value = "'0.0528"
$ws.Range("E44").Value = "  +4.96%  "
$ws.Range("D45").Value = "'19.83"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").Value = "'0.582"
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("D47").Value = "'0.0907"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("D49").Value = "'11.40"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "'16.74"
$ws.Range("E50").Value = "  +3.41%  "
$ws.Range("E51").Value = "  +5.23%  "
